# Season-record columns: Wins / Losses / Ties appended after the last
# existing column (AC), for the single sheet in this workbook.
#
# Header row (row 1) reuses the same bold/bordered/centered header style
# that's already applied to A1:AC1 -- copy A1's format onto the new header
# cells first, then overwrite their text so the shared style index is
# reused instead of a new one being minted.
#
# Data rows (2-47) just get plain numeric values (no special style),
# matching the rest of the data body.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -----------------------------------------------------
$ws.Range("A1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- data rows (2-47) : Wins=79, Losses=83, Ties=0 for every row ----
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 79   # AD
    $ws.Cells.Item($row, 31).Value = 83   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
